# "Fruta / hortaliza, semanal" — weekly price update.
# A new weekly observation is inserted as the new row 3 (pushing the
# previously existing rows 3-6 down to rows 4-7); the new row holds the
# latest price data for the same market/product/variety combination.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 3, shifting rows 3:6 down to 4:7.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with this week's data.
$ws.Cells.Item(3, 1).Value  = 1
$ws.Cells.Item(3, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(3, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(3, 4).Value  = 44487
$ws.Cells.Item(3, 5).Value  = 15
$ws.Cells.Item(3, 6).Value  = "Fruta"
$ws.Cells.Item(3, 7).Value  = 100101
$ws.Cells.Item(3, 8).Value  = "Berries"
$ws.Cells.Item(3, 9).Value  = 100101007
$ws.Cells.Item(3, 10).Value = "Kiwi"
$ws.Cells.Item(3, 11).Value = "Hayward"
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 300
$ws.Cells.Item(3, 14).Value = 14000
$ws.Cells.Item(3, 15).Value = 15000
$ws.Cells.Item(3, 16).Value = 14500
$ws.Cells.Item(3, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(3, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(3, 19).Value = 1450
$ws.Cells.Item(3, 20).Value = 10
